$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 3
Write-Output "Before save: ScrollRow=$($win.ScrollRow) ScrollColumn=$($win.ScrollColumn)"
